# Weekly update: a new week's data row is added for
# "Femacal de La Calera - Ciboulette". The new observation is inserted as
# row 65 (pushing the existing rows 65..212 down to 66..213), and the new
# row is populated with the same values as the row that used to be in that
# spot (now shifted down to row 66), except for a new date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 65
$lastCol = 18   # columns A..R

# Insert a blank row at position 65; everything below (old rows 65-212)
# shifts down by one (to rows 66-213).
$ws.Rows.Item($newRow).Insert()

# The row that used to occupy position 65 is now one row below, at
# $newRow + 1. Duplicate its values into the freshly inserted blank row,
# cell by cell (Range.Value round-trips unreliably in this host, so copy
# via individual cells).
$sourceRow = $newRow + 1
for ($col = 1; $col -le $lastCol; $col++) {
    $ws.Cells.Item($newRow, $col).Value2 = $ws.Cells.Item($sourceRow, $col).Value2
}

# Overwrite the date (column D = 4) of the new row with the new
# observation's date.
$ws.Cells.Item($newRow, 4).Value2 = 44519
